$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# G2 assessment-title value updated from "RMK Nextgen - Add Assessment" to "RMK Nextgen | Add Assessment"
$ws.Range("G2").Value = "RMK Nextgen | Add Assessment"

# New column I header: "add assessment 2", matching the format/style of H1 (page_title/subject header cell)
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "add assessment 2"

# New column I value holds the old G2 text ("RMK Nextgen - Add Assessment"), no special style (like D2)
$ws.Range("I2").Value = "RMK Nextgen - Add Assessment"

# Move the active selection to D15 (matches progress further down the test sheet)
$ws.Range("D15").Select()
